$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----- Header row (row 1) -----
$headers = @(
    "Job_Id",
    "Job_Title",
    "Job_Description",
    "Total_Years_Min_Exp",
    "Total_Years_Max_Exp",
    "Work_Mode",
    "Job_Location",
    "LinkedIn_Poster",
    "LinkedIn_Posted",
    "Resume_received",
    "Resume_downloaded"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Build the header format once on A1 (bold font, thin box border,
# centered/top-aligned) then copy that exact format to the rest of the
# header row so only a single new style entry is produced.
$headerCell = $ws.Range("A1")
$headerCell.Font.Bold = $true
$headerCell.HorizontalAlignment = -4108   # xlCenter
$headerCell.VerticalAlignment = -4160     # xlTop
$headerCell.Borders.LineStyle = 1         # xlContinuous

$headerCell.Copy()
$ws.Range("B1:K1").PasteSpecial(-4122)    # xlPasteFormats
$excel.CutCopyMode = $false

# ----- Data row (row 2) -----
$ws.Range("A2").Value = "JD_001"
$ws.Range("B2").Value = "Junior RPA Developer"
$ws.Range("C2").Value = "We are seeking a Junior RPA Developer to design, develop, and support automation solutions.`nCollaborate with teams to streamline business processes using RPA tools like UiPath or Automation Anywhere. Join Akkodis to grow your skills in a dynamic, tech-driven environment"
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 4
$ws.Range("F2").Value = "Hybrid"
$ws.Range("G2").Value = "Bengaluru, Karnataka, India"

# The multi-line Job_Description triggers an automatic custom row height;
# AutoFit restores the row to its natural (non-custom) height so the
# serialized XML doesn't carry stray ht/customHeight attributes.
$ws.Rows.Item(2).AutoFit() | Out-Null
